$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update as a literal text value (leading apostrophe forces
# Excel to store it as text instead of auto-converting look-alike numbers),
# then ClearFormats() drops the quote-prefix style so the cell keeps the
# workbook-default (unstyled) formatting, matching the source data.
function Set-TextValue($address, $text) {
    $cell = $ws.Range($address)
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

Set-TextValue "D2" "42.129.00"
Set-TextValue "E2" "  +2.91%  "
Set-TextValue "D3" "2.212.52"
Set-TextValue "E3" "  +1.79%  "
Set-TextValue "E4" "  -0.21%  "
Set-TextValue "D5" "251.71"
Set-TextValue "E5" "  +0.44%  "
Set-TextValue "D6" "0.614"
Set-TextValue "E6" "  -0.88%  "
Set-TextValue "D7" "67.37"
Set-TextValue "E7" "  +0.85%  "
Set-TextValue "E8" "  -0.10%  "
Set-TextValue "D9" "0.618"
Set-TextValue "E9" "  +9.44%  "
Set-TextValue "D10" "38.93"
Set-TextValue "E10" "  +5.60%  "
Set-TextValue "D11" "59.36"
Set-TextValue "E11" "  +1.99%  "
Set-TextValue "D12" "0.0936"
Set-TextValue "E12" "  +0.84%  "
Set-TextValue "D13" "7.02"
Set-TextValue "E13" "  +1.20%  "
Set-TextValue "E14" "  +0.14%  "
Set-TextValue "D15" "2.545.39"
Set-TextValue "E15" "  +1.90%  "
Set-TextValue "D16" "0.866"
Set-TextValue "E16" "  +1.05%  "
Set-TextValue "D17" "14.43"
Set-TextValue "E17" "  +0.46%  "
Set-TextValue "D18" "2.200.20"
Set-TextValue "E18" "  +1.26%  "
Set-TextValue "D19" "41.964.24"
Set-TextValue "E19" "  +2.69%  "
Set-TextValue "D20" "0.0₃0960"
Set-TextValue "E20" "  +2.16%  "
Set-TextValue "D21" "72.20"
Set-TextValue "E21" "  +1.20%  "
Set-TextValue "E22" "  -0.58%  "
Set-TextValue "D23" "231.13"
Set-TextValue "E23" "  +0.35%  "
Set-TextValue "E24" "  -2.16%  "
Set-TextValue "E25" "  +0.17%  "
Set-TextValue "E26" "  +0.14%  "
Set-TextValue "D27" "11.12"
Set-TextValue "E27" "  -4.33%  "
Set-TextValue "D28" "2.41"
Set-TextValue "E28" "  -2.60%  "
Set-TextValue "E29" "  -1.07%  "
Set-TextValue "D30" "2.19"
Set-TextValue "E30" "  +1.21%  "
Set-TextValue "D31" "166.81"
Set-TextValue "E32" "  -0.25%  "
Set-TextValue "D33" "5.92"
Set-TextValue "E33" "  +10.44%  "
Set-TextValue "E34" "  +3.95%  "
Set-TextValue "D35" "0.0777"
Set-TextValue "E35" "  +7.80%  "
Set-TextValue "E36" "  +0.64%  "
Set-TextValue "E37" "  +2.23%  "
Set-TextValue "E38" "  +0.81%  "
Set-TextValue "D39" "4.10"
Set-TextValue "E39" "  +2.56%  "
Set-TextValue "D40" "0.0312"
Set-TextValue "E40" "  +5.80%  "
Set-TextValue "E41" "  +1.03%  "
Set-TextValue "B42" "THORChain"
Set-TextValue "C42" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D42" "5.66"
Set-TextValue "E42" "  +0.80%  "
Set-TextValue "B43" "FTXToken"
Set-TextValue "C43" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D43" "5.15"
Set-TextValue "E43" "  +5.48%  "
Set-TextValue "D44" "11.97"
Set-TextValue "D45" "61.18"
Set-TextValue "E45" "  -4.01%  "
Set-TextValue "E46" "  -2.67%  "
Set-TextValue "E47" "  +0.08%  "
Set-TextValue "E48" "  -0.76%  "
Set-TextValue "E49" "  -0.46%  "
Set-TextValue "D50" "1.14"
Set-TextValue "E50" "  +1.50%  "
Set-TextValue "E51" "  +4.31%  "
